$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement every value in columns A and B (rows 2-23) by 1
for ($r = 2; $r -le 23; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $aCell.Value2 = $aCell.Value2 - 1
    $bCell.Value2 = $bCell.Value2 - 1
}

# Update the active selection to E11 (was E21)
$ws.Range("E11").Select()
